$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All touched cells are inline/shared strings in the source (t="inlineStr").
# Force text number format so numeric-looking values (e.g. "0.611", "42.413.62")
# are stored as text rather than being auto-converted to numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.413.62'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.64%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.186.38'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.32%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.11'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.611'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.47%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '75.11'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.45%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.579'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -5.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.85'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0909'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.84%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.73'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.510.93'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.10'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.179.52'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.767'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -5.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.307.10'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.64%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.84'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.84'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '226.54'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.27'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -13.38%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.21%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.41'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -5.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.40'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.77%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.44'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.52%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '171.79'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.58%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.14'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.99'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0818'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.13'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.85%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.106'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.19'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '11.97'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -9.56%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.15'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -8.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.193'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.18%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '58.47'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.23%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.53'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +9.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.62'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0970'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.61%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'WOONetwork'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.459'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.62%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.16'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -4.66%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.51%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.58%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.06%  '
